$wb = $excel.ActiveWorkbook

# Insert a new "is_targeted list" sheet right after "analyte_class list"
# (and before "resolution_x_unit list"), matching the new workbook sheet order.
$afterSheet = $wb.Worksheets.Item("analyte_class list")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "is_targeted list"

# Populate it with the two allowed boolean values.
$newSheet.Range("A1").Value = "TRUE"
$newSheet.Range("A2").Value = "FALSE"

# Point the "is_targeted" column's validation (N2:N1048576) at the new list
# instead of the inline "TRUE,FALSE" formula.
$mainSheet = $wb.Worksheets.Item("Export as TSV")
$targetRange = $mainSheet.Range("N2:N1048576")
$targetRange.Validation.Delete()
$targetRange.Validation.Add(3, 1, 1, "='is_targeted list'!`$A`$1:`$A`$2")
$targetRange.Validation.ErrorTitle = "Value must come from list"
$targetRange.Validation.ErrorMessage = "Value must be one of: TRUE / FALSE."
